$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(114).Insert()

$ws.Cells.Item(114, 1).Value = 4
$ws.Cells.Item(114, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(114, 3).Value = "Los Lagos"
$ws.Cells.Item(114, 4).Value = 45134
$ws.Cells.Item(114, 5).Value = 10
$ws.Cells.Item(114, 6).Value = 100112022
$ws.Cells.Item(114, 7).Value = "Arveja Verde"
$ws.Cells.Item(114, 8).Value = "Perfection"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 40
$ws.Cells.Item(114, 11).Value = 40000
$ws.Cells.Item(114, 12).Value = 40000
$ws.Cells.Item(114, 13).Value = 40000
$ws.Cells.Item(114, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(114, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(114, 16).Value = 1600
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = "Hortaliza"
